$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (AD1:AF1) -------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold font, borders, centered alignment) from
# the existing last header cell (AC1) onto the new header cells without
# touching the text we just wrote.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Season record columns (Wins / Losses / Ties) for every data row -----
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Season record columns (Wins/Losses/Ties) added."
